$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two "Beitragsbemessungsgrenze PV Ost/West" rows to the new
# single "Beitragsbemessungsgrenze PV " / "Jahresarbeitsentgeltgrenze PV "
# labels (both keep a trailing space, as in the source workbook).
$ws.Range("A4").Value = "Beitragsbemessungsgrenze PV "
$ws.Range("A5").Value = "Jahresarbeitsentgeltgrenze PV "

# Update the active selection to A4, as recorded in the saved workbook.
$ws.Range("A4").Select()
